# edit.ps1
# Applies "river update May 2024" changes to Sheet1:
#  1. Updates G/H column values for E. coli / Nitrate rows in the three existing
#     year-range blocks (2016-2020, 2017-2021, 2018-2022).
#  2. Appends a new 16-row block (rows 50-65) for the "2019 - 2023" year range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update recalculated Mean (G) / Annual Maximum (H) values ---

$ws.Cells.Item(4, 7).Value = 606.7196193499429
$ws.Cells.Item(4, 8).Value = 6725.54172879727

$ws.Cells.Item(5, 7).Value = 606.7196193499429
$ws.Cells.Item(5, 8).Value = 6725.54172879727

$ws.Cells.Item(6, 7).Value = 606.7196193499429
$ws.Cells.Item(6, 8).Value = 6725.54172879727

$ws.Cells.Item(7, 7).Value = 606.7196193499429
$ws.Cells.Item(7, 8).Value = 6725.54172879727

$ws.Cells.Item(10, 7).Value = 0.485166568832513

$ws.Cells.Item(11, 7).Value = 0.485166568832513

$ws.Cells.Item(20, 7).Value = 634.8396903185979
$ws.Cells.Item(20, 8).Value = 6725.54172879727

$ws.Cells.Item(21, 7).Value = 634.8396903185979
$ws.Cells.Item(21, 8).Value = 6725.54172879727

$ws.Cells.Item(22, 7).Value = 634.8396903185979
$ws.Cells.Item(22, 8).Value = 6725.54172879727

$ws.Cells.Item(23, 7).Value = 634.8396903185979
$ws.Cells.Item(23, 8).Value = 6725.54172879727

$ws.Cells.Item(26, 7).Value = 0.456165375268079

$ws.Cells.Item(27, 7).Value = 0.456165375268079

$ws.Cells.Item(36, 7).Value = 700.571548898865
$ws.Cells.Item(36, 8).Value = 6725.54172879727

$ws.Cells.Item(37, 7).Value = 700.571548898865
$ws.Cells.Item(37, 8).Value = 6725.54172879727

$ws.Cells.Item(38, 7).Value = 700.571548898865
$ws.Cells.Item(38, 8).Value = 6725.54172879727

$ws.Cells.Item(39, 7).Value = 700.571548898865
$ws.Cells.Item(39, 8).Value = 6725.54172879727

$ws.Cells.Item(42, 7).Value = 0.465995812772702

$ws.Cells.Item(43, 7).Value = 0.465995812772702

# --- Step 2: append new data rows for year range "2019 - 2023" (rows 50-65) ---

# Row 50: DRP (95th Percentile)
$ws.Cells.Item(50, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(50, 2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(50, 3).Value = "D"
$ws.Cells.Item(50, 4).Value = "2019 - 2023"
$ws.Cells.Item(50, 5).Value = "RepSite"
$ws.Cells.Item(50, 6).Value = 0.08
$ws.Cells.Item(50, 7).Value = 0.0872037037037037
$ws.Cells.Item(50, 8).Value = 0.173
$ws.Cells.Item(50, 9).Value = 0.155
$ws.Cells.Item(50, 12).Value = 0.102
$ws.Cells.Item(50, 13).Value = 0.11924
$ws.Cells.Item(50, 14).Value = 0.13418
$ws.Cells.Item(50, 15).Value = 1789261
$ws.Cells.Item(50, 16).Value = 5528869
$ws.Cells.Item(50, 17).Value = "Manawatu District"
$ws.Cells.Item(50, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(50, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(50, 20).Value = "West_6"
$ws.Cells.Item(50, 21).Value = "mg/L"

# Row 51: DRP (Median)
$ws.Cells.Item(51, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(51, 2).Value = "DRP (Median)"
$ws.Cells.Item(51, 3).Value = "D"
$ws.Cells.Item(51, 4).Value = "2019 - 2023"
$ws.Cells.Item(51, 5).Value = "RepSite"
$ws.Cells.Item(51, 6).Value = 0.08
$ws.Cells.Item(51, 7).Value = 0.0872037037037037
$ws.Cells.Item(51, 8).Value = 0.173
$ws.Cells.Item(51, 9).Value = 0.155
$ws.Cells.Item(51, 12).Value = 0.102
$ws.Cells.Item(51, 13).Value = 0.11924
$ws.Cells.Item(51, 14).Value = 0.13418
$ws.Cells.Item(51, 15).Value = 1789261
$ws.Cells.Item(51, 16).Value = 5528869
$ws.Cells.Item(51, 17).Value = "Manawatu District"
$ws.Cells.Item(51, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(51, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(51, 20).Value = "West_6"
$ws.Cells.Item(51, 21).Value = "mg/L"

# Row 52: E coli (>260)
$ws.Cells.Item(52, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(52, 2).Value = "E coli (>260)"
$ws.Cells.Item(52, 3).Value = "E"
$ws.Cells.Item(52, 4).Value = "2019 - 2023"
$ws.Cells.Item(52, 5).Value = "RepSite"
$ws.Cells.Item(52, 6).Value = 335
$ws.Cells.Item(52, 7).Value = 719.502774743225
$ws.Cells.Item(52, 8).Value = 6725.54172879727
$ws.Cells.Item(52, 9).Value = 2943.6
$ws.Cells.Item(52, 10).Value = 33.3333333333333
$ws.Cells.Item(52, 11).Value = 61.1111111111111
$ws.Cells.Item(52, 12).Value = 96
$ws.Cells.Item(52, 13).Value = 1189.28
$ws.Cells.Item(52, 14).Value = 2221.24
$ws.Cells.Item(52, 15).Value = 1789261
$ws.Cells.Item(52, 16).Value = 5528869
$ws.Cells.Item(52, 17).Value = "Manawatu District"
$ws.Cells.Item(52, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(52, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(52, 20).Value = "West_6"
$ws.Cells.Item(52, 21).Value = "% exceedances over 260/100 mL"

# Row 53: E coli (>540)
$ws.Cells.Item(53, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(53, 2).Value = "E coli (>540)"
$ws.Cells.Item(53, 3).Value = "E"
$ws.Cells.Item(53, 4).Value = "2019 - 2023"
$ws.Cells.Item(53, 5).Value = "RepSite"
$ws.Cells.Item(53, 6).Value = 335
$ws.Cells.Item(53, 7).Value = 719.502774743225
$ws.Cells.Item(53, 8).Value = 6725.54172879727
$ws.Cells.Item(53, 9).Value = 2943.6
$ws.Cells.Item(53, 10).Value = 33.3333333333333
$ws.Cells.Item(53, 11).Value = 61.1111111111111
$ws.Cells.Item(53, 12).Value = 96
$ws.Cells.Item(53, 13).Value = 1189.28
$ws.Cells.Item(53, 14).Value = 2221.24
$ws.Cells.Item(53, 15).Value = 1789261
$ws.Cells.Item(53, 16).Value = 5528869
$ws.Cells.Item(53, 17).Value = "Manawatu District"
$ws.Cells.Item(53, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(53, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(53, 20).Value = "West_6"
$ws.Cells.Item(53, 21).Value = "% exceedances over 540/100 mL"

# Row 54: E coli (Median)
$ws.Cells.Item(54, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(54, 2).Value = "E coli (Median)"
$ws.Cells.Item(54, 3).Value = "E"
$ws.Cells.Item(54, 4).Value = "2019 - 2023"
$ws.Cells.Item(54, 5).Value = "RepSite"
$ws.Cells.Item(54, 6).Value = 335
$ws.Cells.Item(54, 7).Value = 719.502774743225
$ws.Cells.Item(54, 8).Value = 6725.54172879727
$ws.Cells.Item(54, 9).Value = 2943.6
$ws.Cells.Item(54, 10).Value = 33.3333333333333
$ws.Cells.Item(54, 11).Value = 61.1111111111111
$ws.Cells.Item(54, 12).Value = 96
$ws.Cells.Item(54, 13).Value = 1189.28
$ws.Cells.Item(54, 14).Value = 2221.24
$ws.Cells.Item(54, 15).Value = 1789261
$ws.Cells.Item(54, 16).Value = 5528869
$ws.Cells.Item(54, 17).Value = "Manawatu District"
$ws.Cells.Item(54, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(54, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(54, 20).Value = "West_6"
$ws.Cells.Item(54, 21).Value = "E. coli/100 mL"

# Row 55: E coli (95th Percentile)
$ws.Cells.Item(55, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(55, 2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(55, 3).Value = "E"
$ws.Cells.Item(55, 4).Value = "2019 - 2023"
$ws.Cells.Item(55, 5).Value = "RepSite"
$ws.Cells.Item(55, 6).Value = 335
$ws.Cells.Item(55, 7).Value = 719.502774743225
$ws.Cells.Item(55, 8).Value = 6725.54172879727
$ws.Cells.Item(55, 9).Value = 2943.6
$ws.Cells.Item(55, 10).Value = 33.3333333333333
$ws.Cells.Item(55, 11).Value = 61.1111111111111
$ws.Cells.Item(55, 12).Value = 96
$ws.Cells.Item(55, 13).Value = 1189.28
$ws.Cells.Item(55, 14).Value = 2221.24
$ws.Cells.Item(55, 15).Value = 1789261
$ws.Cells.Item(55, 16).Value = 5528869
$ws.Cells.Item(55, 17).Value = "Manawatu District"
$ws.Cells.Item(55, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(55, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(55, 20).Value = "West_6"
$ws.Cells.Item(55, 21).Value = "E. coli/100 mL"

# Row 56: Ammoniacal-N (95th Percentile)
$ws.Cells.Item(56, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(56, 2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(56, 3).Value = "C"
$ws.Cells.Item(56, 4).Value = "2019 - 2023"
$ws.Cells.Item(56, 5).Value = "RepSite"
$ws.Cells.Item(56, 6).Value = 0.2191
$ws.Cells.Item(56, 7).Value = 0.268325743806501
$ws.Cells.Item(56, 8).Value = 1.46306340177539
$ws.Cells.Item(56, 9).Value = 0.73445
$ws.Cells.Item(56, 12).Value = 0.29699
$ws.Cells.Item(56, 13).Value = 0.34885
$ws.Cells.Item(56, 14).Value = 0.54176
$ws.Cells.Item(56, 15).Value = 1789261
$ws.Cells.Item(56, 16).Value = 5528869
$ws.Cells.Item(56, 17).Value = "Manawatu District"
$ws.Cells.Item(56, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(56, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(56, 20).Value = "West_6"
$ws.Cells.Item(56, 21).Value = "mg NH4-N/L"

# Row 57: Ammoniacal-N (Median)
$ws.Cells.Item(57, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(57, 2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(57, 3).Value = "B"
$ws.Cells.Item(57, 4).Value = "2019 - 2023"
$ws.Cells.Item(57, 5).Value = "RepSite"
$ws.Cells.Item(57, 6).Value = 0.2191
$ws.Cells.Item(57, 7).Value = 0.268325743806501
$ws.Cells.Item(57, 8).Value = 1.46306340177539
$ws.Cells.Item(57, 9).Value = 0.73445
$ws.Cells.Item(57, 12).Value = 0.29699
$ws.Cells.Item(57, 13).Value = 0.34885
$ws.Cells.Item(57, 14).Value = 0.54176
$ws.Cells.Item(57, 15).Value = 1789261
$ws.Cells.Item(57, 16).Value = 5528869
$ws.Cells.Item(57, 17).Value = "Manawatu District"
$ws.Cells.Item(57, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(57, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(57, 20).Value = "West_6"
$ws.Cells.Item(57, 21).Value = "mg NH4-N/L"

# Row 58: Nitrate-N (95th Percentile)
$ws.Cells.Item(58, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(58, 2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(58, 3).Value = "A"
$ws.Cells.Item(58, 4).Value = "2019 - 2023"
$ws.Cells.Item(58, 5).Value = "RepSite"
$ws.Cells.Item(58, 6).Value = 0.2295
$ws.Cells.Item(58, 7).Value = 0.46927523536412
$ws.Cells.Item(58, 8).Value = 1.52
$ws.Cells.Item(58, 9).Value = 1.386
$ws.Cells.Item(58, 12).Value = 0.02773
$ws.Cells.Item(58, 13).Value = 1.1232
$ws.Cells.Item(58, 14).Value = 1.2798
$ws.Cells.Item(58, 15).Value = 1789261
$ws.Cells.Item(58, 16).Value = 5528869
$ws.Cells.Item(58, 17).Value = "Manawatu District"
$ws.Cells.Item(58, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(58, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(58, 20).Value = "West_6"
$ws.Cells.Item(58, 21).Value = "mg NO3-N/L"

# Row 59: Nitrate-N (Median)
$ws.Cells.Item(59, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(59, 2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(59, 3).Value = "A"
$ws.Cells.Item(59, 4).Value = "2019 - 2023"
$ws.Cells.Item(59, 5).Value = "RepSite"
$ws.Cells.Item(59, 6).Value = 0.2295
$ws.Cells.Item(59, 7).Value = 0.46927523536412
$ws.Cells.Item(59, 8).Value = 1.52
$ws.Cells.Item(59, 9).Value = 1.386
$ws.Cells.Item(59, 12).Value = 0.02773
$ws.Cells.Item(59, 13).Value = 1.1232
$ws.Cells.Item(59, 14).Value = 1.2798
$ws.Cells.Item(59, 15).Value = 1789261
$ws.Cells.Item(59, 16).Value = 5528869
$ws.Cells.Item(59, 17).Value = "Manawatu District"
$ws.Cells.Item(59, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(59, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(59, 20).Value = "West_6"
$ws.Cells.Item(59, 21).Value = "mg NO3-N/L"

# Row 60: Soluble Inorganic Nitrogen (95th Percentile)
$ws.Cells.Item(60, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(60, 2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(60, 4).Value = "2019 - 2023"
$ws.Cells.Item(60, 5).Value = "RepSite"
$ws.Cells.Item(60, 6).Value = 0.78
$ws.Cells.Item(60, 7).Value = 0.917777777777778
$ws.Cells.Item(60, 8).Value = 1.9
$ws.Cells.Item(60, 9).Value = 1.79
$ws.Cells.Item(60, 12).Value = 0.65
$ws.Cells.Item(60, 13).Value = 1.5328
$ws.Cells.Item(60, 14).Value = 1.7818
$ws.Cells.Item(60, 15).Value = 1789261
$ws.Cells.Item(60, 16).Value = 5528869
$ws.Cells.Item(60, 17).Value = "Manawatu District"
$ws.Cells.Item(60, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(60, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(60, 20).Value = "West_6"
$ws.Cells.Item(60, 21).Value = "g/m3"

# Row 61: Soluble Inorganic Nitrogen (Median)
$ws.Cells.Item(61, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(61, 2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(61, 4).Value = "2019 - 2023"
$ws.Cells.Item(61, 5).Value = "RepSite"
$ws.Cells.Item(61, 6).Value = 0.78
$ws.Cells.Item(61, 7).Value = 0.917777777777778
$ws.Cells.Item(61, 8).Value = 1.9
$ws.Cells.Item(61, 9).Value = 1.79
$ws.Cells.Item(61, 12).Value = 0.65
$ws.Cells.Item(61, 13).Value = 1.5328
$ws.Cells.Item(61, 14).Value = 1.7818
$ws.Cells.Item(61, 15).Value = 1789261
$ws.Cells.Item(61, 16).Value = 5528869
$ws.Cells.Item(61, 17).Value = "Manawatu District"
$ws.Cells.Item(61, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(61, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(61, 20).Value = "West_6"
$ws.Cells.Item(61, 21).Value = "g/m3"

# Row 62: Total Nitrogen (95th Percentile)
$ws.Cells.Item(62, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(62, 2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(62, 4).Value = "2019 - 2023"
$ws.Cells.Item(62, 5).Value = "RepSite"
$ws.Cells.Item(62, 6).Value = 2.105
$ws.Cells.Item(62, 7).Value = 2.23592592592593
$ws.Cells.Item(62, 8).Value = 7.77
$ws.Cells.Item(62, 9).Value = 3.134
$ws.Cells.Item(62, 12).Value = 1.91
$ws.Cells.Item(62, 13).Value = 2.826
$ws.Cells.Item(62, 14).Value = 3.0454
$ws.Cells.Item(62, 15).Value = 1789261
$ws.Cells.Item(62, 16).Value = 5528869
$ws.Cells.Item(62, 17).Value = "Manawatu District"
$ws.Cells.Item(62, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(62, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(62, 20).Value = "West_6"
$ws.Cells.Item(62, 21).Value = "g/m3"

# Row 63: Total Nitrogen (Median)
$ws.Cells.Item(63, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(63, 2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(63, 4).Value = "2019 - 2023"
$ws.Cells.Item(63, 5).Value = "RepSite"
$ws.Cells.Item(63, 6).Value = 2.105
$ws.Cells.Item(63, 7).Value = 2.23592592592593
$ws.Cells.Item(63, 8).Value = 7.77
$ws.Cells.Item(63, 9).Value = 3.134
$ws.Cells.Item(63, 12).Value = 1.91
$ws.Cells.Item(63, 13).Value = 2.826
$ws.Cells.Item(63, 14).Value = 3.0454
$ws.Cells.Item(63, 15).Value = 1789261
$ws.Cells.Item(63, 16).Value = 5528869
$ws.Cells.Item(63, 17).Value = "Manawatu District"
$ws.Cells.Item(63, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(63, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(63, 20).Value = "West_6"
$ws.Cells.Item(63, 21).Value = "g/m3"

# Row 64: Total Phosphorus (95th Percentile)
$ws.Cells.Item(64, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(64, 2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(64, 4).Value = "2019 - 2023"
$ws.Cells.Item(64, 5).Value = "RepSite"
$ws.Cells.Item(64, 6).Value = 0.246
$ws.Cells.Item(64, 7).Value = 0.318592592592593
$ws.Cells.Item(64, 8).Value = 0.984
$ws.Cells.Item(64, 9).Value = 0.7158
$ws.Cells.Item(64, 12).Value = 0.342
$ws.Cells.Item(64, 13).Value = 0.41292
$ws.Cells.Item(64, 14).Value = 0.6509200000000001
$ws.Cells.Item(64, 15).Value = 1789261
$ws.Cells.Item(64, 16).Value = 5528869
$ws.Cells.Item(64, 17).Value = "Manawatu District"
$ws.Cells.Item(64, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(64, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(64, 20).Value = "West_6"
$ws.Cells.Item(64, 21).Value = "g/m3"

# Row 65: Total Phosphorus (Median)
$ws.Cells.Item(65, 1).Value = "Kaikokopu Stream at Himatangi Beach"
$ws.Cells.Item(65, 2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(65, 4).Value = "2019 - 2023"
$ws.Cells.Item(65, 5).Value = "RepSite"
$ws.Cells.Item(65, 6).Value = 0.246
$ws.Cells.Item(65, 7).Value = 0.318592592592593
$ws.Cells.Item(65, 8).Value = 0.984
$ws.Cells.Item(65, 9).Value = 0.7158
$ws.Cells.Item(65, 12).Value = 0.342
$ws.Cells.Item(65, 13).Value = 0.41292
$ws.Cells.Item(65, 14).Value = 0.6509200000000001
$ws.Cells.Item(65, 15).Value = 1789261
$ws.Cells.Item(65, 16).Value = 5528869
$ws.Cells.Item(65, 17).Value = "Manawatu District"
$ws.Cells.Item(65, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(65, 19).Value = "Northern Manawatu Lakes"
$ws.Cells.Item(65, 20).Value = "West_6"
$ws.Cells.Item(65, 21).Value = "g/m3"

